$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final FIM March 2023: update "current" (rows 7, 8, 16) and
# "difference" (rows 35, 36, 44) blocks for columns J:R (2023 Q1 .. 2025 Q1)

# Row 7: Federal Health Outlays Contribution - current
$ws.Range("J7").Value = -0.0784
$ws.Range("K7").Value = 0.0017
$ws.Range("L7").Value = -0.0419
$ws.Range("M7").Value = -0.0414
$ws.Range("N7").Value = -0.0628
$ws.Range("O7").Value = -0.024
$ws.Range("P7").Value = 0.0267
$ws.Range("Q7").Value = 0.0308
$ws.Range("R7").Value = 0.041

# Row 8: Federal Non Corporate Taxes Contribution - current
$ws.Range("J8").Value = -0.1448
$ws.Range("K8").Value = 0.0079
$ws.Range("L8").Value = 0.0625
$ws.Range("M8").Value = 0.1251
$ws.Range("N8").Value = 0.4023
$ws.Range("O8").Value = 0.3972
$ws.Range("P8").Value = 0.4294
$ws.Range("Q8").Value = 0.3823
$ws.Range("R8").Value = 0.2794

# Row 16: Fiscal Impact - current
$ws.Range("J16").Value = -1.4594
$ws.Range("K16").Value = -1.5607
$ws.Range("L16").Value = -0.4575
$ws.Range("M16").Value = -0.4824
$ws.Range("N16").Value = -0.5423
$ws.Range("O16").Value = -0.1539
$ws.Range("P16").Value = -0.1131
$ws.Range("Q16").Value = 0.0286
$ws.Range("R16").Value = 0.0064

# Row 35: Federal Health Outlays Contribution - difference
$ws.Range("J35").Value = 0.031
$ws.Range("K35").Value = 0.0709
$ws.Range("L35").Value = 0.0956
$ws.Range("M35").Value = 0.0972
$ws.Range("N35").Value = 0.0671
$ws.Range("O35").Value = 0.0294
$ws.Range("P35").Value = 0.0068
$ws.Range("Q35").Value = 0.0068
$ws.Range("R35").Value = 0.0069

# Row 36: Federal Non Corporate Taxes Contribution - difference
$ws.Range("J36").Value = -0.0222
$ws.Range("K36").Value = -0.0331
$ws.Range("L36").Value = -0.0307
$ws.Range("M36").Value = -0.0255
$ws.Range("N36").Value = -0.0203
$ws.Range("O36").Value = -0.0202
$ws.Range("P36").Value = -0.02
$ws.Range("Q36").Value = -0.0209
$ws.Range("R36").Value = -0.0093

# Row 44: Fiscal Impact - difference
$ws.Range("J44").Value = 0.0287
$ws.Range("K44").Value = 0.0454
$ws.Range("L44").Value = 0.0821
$ws.Range("M44").Value = 0.08
$ws.Range("N44").Value = 0.0545
$ws.Range("O44").Value = 0.0236
$ws.Range("P44").Value = -0.0076
$ws.Range("Q44").Value = -0.0143
$ws.Range("R44").Value = -0.0025
